$wb = $excel.ActiveWorkbook
$wsRequest = $wb.Worksheets.Item("Request")

# The "Request" table (Table16) gains a new row (RequestID 8 / ESTIMATE /
# "request an estimation among multiple Synopses") - resize the table so it
# covers the new row, matching what Excel does when a table grows.
$lo = $wsRequest.ListObjects.Item(1)
$lo.Resize($wsRequest.Range("A1:C9"))

# Match the formatting of the existing OperationType column cells (copy B8's
# style onto the freshly added B9) before filling in the new row's values.
$wsRequest.Range("B8").Copy()
$wsRequest.Range("B9").PasteSpecial(-4122)

$wsRequest.Range("A9").Value = 8
$wsRequest.Range("B9").Value = "ESTIMATE"
$wsRequest.Range("C9").Value = "request an estimation among multiple Synopses"

# Switch focus to the "Request" sheet and select the next empty row, A10 -
# reflecting where the user continued working.
$wsRequest.Activate()
$wsRequest.Range("A10").Select()
